# Atualização de bases das ligas, do dia: 07-04-2024 às 22:30
# Swap the contents of columns B:AC between row 11 <-> row 12
# and between row 83 <-> row 84 (column A, the running "id", stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: ".Value" getter is unreliable in this COM-interop runtime (it
# returns a placeholder string instead of the real data). ".Value2"
# works correctly for both get and set, so use that instead.
#
# NOTE: wrapping this logic in a PowerShell function caused the swap to
# only take effect on the first column of the range in this runtime, so
# the swap is inlined (duplicated) for each row pair instead of using a
# helper function.

$rangeA = $ws.Range("B11:AC11")
$rangeB = $ws.Range("B12:AC12")
$valuesA = $rangeA.Value2
$valuesB = $rangeB.Value2
$rangeA.Value2 = $valuesB
$rangeB.Value2 = $valuesA

$rangeC = $ws.Range("B83:AC83")
$rangeD = $ws.Range("B84:AC84")
$valuesC = $rangeC.Value2
$valuesD = $rangeD.Value2
$rangeC.Value2 = $valuesD
$rangeD.Value2 = $valuesC
